$d = $word.ActiveDocument

# 1. Remove the existing (empty) "_GoBack" bookmark that currently sits
#    between ". " and "The HMS shall " in the "daily schedules" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append new content to the paragraph that ends with
#    "Staff shall have access to all schedules maintained by the HMS."
#    Final content: "...maintained by the HMS. Housekeeping staff shall be
#    able to use HMS to manage room status." where the new sentence (not
#    the leading space) is wrapped in a (new) "_GoBack" bookmark.
$oldText = "Staff shall have access to all schedules maintained by the HMS."
$newSentence = "Housekeeping staff shall be able to use HMS to manage room status."
$newText = $oldText + " " + $newSentence

$range = $d.Content
$range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# 3. Re-find the newly inserted sentence and wrap it with the "_GoBack" bookmark.
$bmRange = $d.Content
$bmRange.Find.Execute($newSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
